$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before F to make room for the new question column.
$ws.Columns("F:F").Insert()

# 2. Set the new column F width (matches target 26.0703125)
$ws.Columns("F:F").ColumnWidth = 26.0703125

# 3. Fill in the new header cells for column F (row2 Chinese / row3 English)
$ws.Range("F2").Value = "以何种方式知道iDVC/FreeDIC?"
$ws.Range("F3").Value = "论文，网页检索、推荐、其它"

# 4. Clear the old "填写黑色字体部分" note cell (now at D1, unaffected by the column insert since it's before column F)
$ws.Range("D1").ClearContents()

